$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "image" header in H1, matching the style of the existing headers (e.g. G1)
$ws.Range("H1").Value = "image"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Per-row updates: refreshed follower counts (col F), a couple of popularity (col E)
# tweaks, and the new Spotify image URL (col H).
$rows = @(
    @{ Row = 2; E = $null; F = 12264838; H = "https://i.scdn.co/image/ab6761610000e5ebee452efcf24aa4124fb28d94" },
    @{ Row = 3; E = $null; F = 15114832; H = "https://i.scdn.co/image/ab6761610000e5eb30122c0d3ead72f96bb5ee93" },
    @{ Row = 4; E = $null; F = 8144120; H = "https://i.scdn.co/image/ab6761610000e5eb504ff11d788162fbf8078654" },
    @{ Row = 5; E = $null; F = 13928327; H = "https://i.scdn.co/image/ab6761610000e5eb24e41f491b129093a6fee383" },
    @{ Row = 6; E = $null; F = 8340304; H = "https://i.scdn.co/image/ab6761610000e5eb547d2b41c9f2c97318aad0ed" },
    @{ Row = 7; E = $null; F = 5882877; H = "https://i.scdn.co/image/ab6761610000e5eba6ab3c4df02cec59758ae3fa" },
    @{ Row = 8; E = $null; F = 13731066; H = "https://i.scdn.co/image/ab6761610000e5eb35ca7d2181258b51c0f2cf9e" },
    @{ Row = 9; E = 85; F = 13575711; H = "https://i.scdn.co/image/ab6761610000e5ebc63aded6f4bf4d06d1377106" },
    @{ Row = 10; E = $null; F = 75467555; H = "https://i.scdn.co/image/ab6761610000e5eb4293385d324db8558179afd9" },
    @{ Row = 11; E = $null; F = 22479986; H = "https://i.scdn.co/image/ab6761610000e5ebe707b87e3f65997f6c09bfff" },
    @{ Row = 12; E = $null; F = 29189277; H = "https://i.scdn.co/image/ab6761610000e5eb1908e1a8b79abf71d5598944" },
    @{ Row = 13; E = $null; F = 20194806; H = "https://i.scdn.co/image/ab6761610000e5ebadd503b411a712e277895c8a" },
    @{ Row = 14; E = $null; F = 14188438; H = "https://i.scdn.co/image/ab6761610000e5eb6cad3eff5adc29e20f189a6c" },
    @{ Row = 15; E = $null; F = 24059100; H = "https://i.scdn.co/image/ab6761610000e5eb437b9e2a82505b3d93ff1022" },
    @{ Row = 16; E = $null; F = 9363182; H = "https://i.scdn.co/image/ab6761610000e5eb9c30c6b69a55d48decd71600" },
    @{ Row = 17; E = $null; F = 7234151; H = "https://i.scdn.co/image/ab6761610000e5eb238b2a30c741d42a4c91b7b7" },
    @{ Row = 18; E = $null; F = 20137321; H = "https://i.scdn.co/image/ab6761610000e5eb867008a971fae0f4d913f63a" },
    @{ Row = 19; E = $null; F = 8376711; H = "https://i.scdn.co/image/ab6761610000e5ebc75afcd5a9027f60eaebb5e4" },
    @{ Row = 20; E = $null; F = 1281076; H = "https://i.scdn.co/image/ab6761610000e5eb1ff1685224034e6c12538722" },
    @{ Row = 21; E = $null; F = 71020309; H = "https://i.scdn.co/image/ab6761610000e5eba00b11c129b27a88fc72f36b" },
    @{ Row = 22; E = $null; F = 12504534; H = "https://i.scdn.co/image/ab6761610000e5eb597f9edd2cd1a892d4412b09" },
    @{ Row = 23; E = 84; F = 5589335; H = "https://i.scdn.co/image/ab6761610000e5eba0461c1f2218374aa672ce4e" },
    @{ Row = 24; E = $null; F = 5272350; H = "https://i.scdn.co/image/ab6761610000e5eb803f228472451496cb2f5b88" },
    @{ Row = 25; E = $null; F = 8690894; H = "https://i.scdn.co/image/ab6761610000e5eb32f8f4df5e7b05a7e4d170ca" },
    @{ Row = 26; E = $null; F = 6488601; H = "https://i.scdn.co/image/ab6761610000e5eb876faa285687786c3d314ae0" },
    @{ Row = 27; E = $null; F = 2767861; H = "https://i.scdn.co/image/ab6761610000e5ebc5b88a3924d8318f25f20594" },
    @{ Row = 28; E = $null; F = 2013768; H = "https://i.scdn.co/image/ab6761610000e5eb5e93db92ca7864585fbe5f28" },
    @{ Row = 29; E = $null; F = 11548636; H = "https://i.scdn.co/image/ab6761610000e5eb8278b782cbb5a3963db88ada" },
    @{ Row = 30; E = $null; F = 10818217; H = "https://i.scdn.co/image/ab6761610000e5eb89dccc0a0b3d818c8d5fb711" },
    @{ Row = 31; E = $null; F = 2462095; H = "https://i.scdn.co/image/ab6761610000e5eba36d13951ab8280a63bac16c" },
    @{ Row = 32; E = $null; F = 2892760; H = "https://i.scdn.co/image/ab6761610000e5ebeb63bf6379a9ea8453a30020" },
    @{ Row = 33; E = $null; F = 13184018; H = "https://i.scdn.co/image/ab6761610000e5ebf4593f7b778219838d858c34" },
    @{ Row = 34; E = $null; F = 764444; H = "https://i.scdn.co/image/ab6761610000e5ebaef8392a62c123944b6383b2" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($null -ne $r.E) {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    }
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 8).Value = $r.H
}
